$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.012617333333333
$ws.Cells.Item(2, 8).Value = 3.037852
$ws.Cells.Item(2, 9).Value = 0.0186050446061446
$ws.Cells.Item(2, 10).Value = 0.0186050446061446
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.713818333333333
$ws.Cells.Item(2, 14).Value = 8.141454999999999
$ws.Cells.Item(2, 15).Value = 0.4320628550779991
$ws.Cells.Item(2, 16).Value = 0.4320628550779991
$ws.Cells.Item(2, 17).Value = 2.74805948385111
$ws.Cells.Item(2, 18).Value = 24.73253535466
$ws.Cells.Item(2, 19).Value = 0.008038548691384364
$ws.Cells.Item(2, 20).Value = 0.008038548691384362
$ws.Cells.Item(3, 7).Value = 1.012617333333333
$ws.Cells.Item(3, 8).Value = 3.037852
$ws.Cells.Item(3, 9).Value = 0.0186050446061446
$ws.Cells.Item(3, 10).Value = 0.0186050446061446
$ws.Cells.Item(3, 15).Value = 0.3580555954903459
$ws.Cells.Item(3, 16).Value = 0.3580555954903459
$ws.Cells.Item(3, 17).Value = 2.277349379537778
$ws.Cells.Item(3, 18).Value = 20.49614441584
$ws.Cells.Item(3, 19).Value = 0.006661640325577554
$ws.Cells.Item(3, 20).Value = 0.006661640325577553
$ws.Cells.Item(4, 7).Value = 1.012617333333333
$ws.Cells.Item(4, 8).Value = 3.037852
$ws.Cells.Item(4, 9).Value = 0.0186050446061446
$ws.Cells.Item(4, 10).Value = 0.0186050446061446
$ws.Cells.Item(4, 13).Value = 1.308511666666667
$ws.Cells.Item(4, 14).Value = 3.925535
$ws.Cells.Item(4, 15).Value = 0.2083261357839125
$ws.Cells.Item(4, 16).Value = 0.2083261357839125
$ws.Cells.Item(4, 17).Value = 1.325021594535555
$ws.Cells.Item(4, 18).Value = 11.92519435082
$ws.Cells.Item(4, 19).Value = 0.00387591704888543
$ws.Cells.Item(4, 20).Value = 0.003875917048885429
$ws.Cells.Item(5, 7).Value = 1.012617333333333
$ws.Cells.Item(5, 8).Value = 3.037852
$ws.Cells.Item(5, 9).Value = 0.0186050446061446
$ws.Cells.Item(5, 10).Value = 0.0186050446061446
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.009769666666666668
$ws.Cells.Item(5, 14).Value = 0.029309
$ws.Cells.Item(5, 15).Value = 0.001555413647742459
$ws.Cells.Item(5, 16).Value = 0.001555413647742458
$ws.Cells.Item(5, 17).Value = 0.009892933807555557
$ws.Cells.Item(5, 18).Value = 0.08903640426800001
$ws.Cells.Item(5, 19).Value = 0.00002893854029725453
$ws.Cells.Item(5, 20).Value = 0.00002893854029725452
$ws.Cells.Item(6, 9).Value = 0.7824865355506074
$ws.Cells.Item(6, 10).Value = 0.7824865355506075
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 2.713818333333333
$ws.Cells.Item(6, 14).Value = 8.141454999999999
$ws.Cells.Item(6, 15).Value = 0.4320628550779991
$ws.Cells.Item(6, 16).Value = 0.4320628550779991
$ws.Cells.Item(6, 17).Value = 115.5772313652755
$ws.Cells.Item(6, 18).Value = 1040.19508228748
$ws.Cells.Item(6, 19).Value = 0.3380833666100877
$ws.Cells.Item(6, 20).Value = 0.3380833666100877
$ws.Cells.Item(7, 9).Value = 0.7824865355506074
$ws.Cells.Item(7, 10).Value = 0.7824865355506075
$ws.Cells.Item(7, 15).Value = 0.3580555954903459
$ws.Cells.Item(7, 16).Value = 0.3580555954903459
$ws.Cells.Item(7, 19).Value = 0.2801736824497505
$ws.Cells.Item(7, 20).Value = 0.2801736824497505
$ws.Cells.Item(8, 9).Value = 0.7824865355506074
$ws.Cells.Item(8, 10).Value = 0.7824865355506075
$ws.Cells.Item(8, 13).Value = 1.308511666666667
$ws.Cells.Item(8, 14).Value = 3.925535
$ws.Cells.Item(8, 15).Value = 0.2083261357839125
$ws.Cells.Item(8, 16).Value = 0.2083261357839125
$ws.Cells.Item(8, 17).Value = 55.72744269021778
$ws.Cells.Item(8, 18).Value = 501.54698421196
$ws.Cells.Item(8, 19).Value = 0.1630123962541992
$ws.Cells.Item(8, 20).Value = 0.1630123962541991
$ws.Cells.Item(9, 9).Value = 0.7824865355506074
$ws.Cells.Item(9, 10).Value = 0.7824865355506075
$ws.Cells.Item(9, 11).Value = 1
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.009769666666666668
$ws.Cells.Item(9, 14).Value = 0.029309
$ws.Cells.Item(9, 15).Value = 0.001555413647742459
$ws.Cells.Item(9, 16).Value = 0.001555413647742458
$ws.Cells.Item(9, 17).Value = 0.4160746542337778
$ws.Cells.Item(9, 18).Value = 3.744671888104
$ws.Cells.Item(9, 19).Value = 0.001217090236570129
$ws.Cells.Item(9, 20).Value = 0.001217090236570129
$ws.Cells.Item(10, 7).Value = 10.82599466666667
$ws.Cells.Item(10, 8).Value = 32.477984
$ws.Cells.Item(10, 9).Value = 0.198908419843248
$ws.Cells.Item(10, 10).Value = 0.198908419843248
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 2.713818333333333
$ws.Cells.Item(10, 14).Value = 8.141454999999999
$ws.Cells.Item(10, 15).Value = 0.4320628550779991
$ws.Cells.Item(10, 16).Value = 0.4320628550779991
$ws.Cells.Item(10, 17).Value = 29.37978280296888
$ws.Cells.Item(10, 18).Value = 264.41804522672
$ws.Cells.Item(10, 19).Value = 0.08594093977652707
$ws.Cells.Item(10, 20).Value = 0.08594093977652706
$ws.Cells.Item(11, 7).Value = 10.82599466666667
$ws.Cells.Item(11, 8).Value = 32.477984
$ws.Cells.Item(11, 9).Value = 0.198908419843248
$ws.Cells.Item(11, 10).Value = 0.198908419843248
$ws.Cells.Item(11, 15).Value = 0.3580555954903459
$ws.Cells.Item(11, 16).Value = 0.3580555954903459
$ws.Cells.Item(11, 17).Value = 24.34737331214222
$ws.Cells.Item(11, 18).Value = 219.12635980928
$ws.Cells.Item(11, 19).Value = 0.07122027271501791
$ws.Cells.Item(11, 20).Value = 0.07122027271501789
$ws.Cells.Item(12, 7).Value = 10.82599466666667
$ws.Cells.Item(12, 8).Value = 32.477984
$ws.Cells.Item(12, 9).Value = 0.198908419843248
$ws.Cells.Item(12, 10).Value = 0.198908419843248
$ws.Cells.Item(12, 13).Value = 1.308511666666667
$ws.Cells.Item(12, 14).Value = 3.925535
$ws.Cells.Item(12, 15).Value = 0.2083261357839125
$ws.Cells.Item(12, 16).Value = 0.2083261357839125
$ws.Cells.Item(12, 17).Value = 14.16594032460445
$ws.Cells.Item(12, 18).Value = 127.49346292144
$ws.Cells.Item(12, 19).Value = 0.04143782248082797
$ws.Cells.Item(12, 20).Value = 0.04143782248082795
$ws.Cells.Item(13, 7).Value = 10.82599466666667
$ws.Cells.Item(13, 8).Value = 32.477984
$ws.Cells.Item(13, 9).Value = 0.198908419843248
$ws.Cells.Item(13, 10).Value = 0.198908419843248
$ws.Cells.Item(13, 11).Value = 1
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.009769666666666668
$ws.Cells.Item(13, 14).Value = 0.029309
$ws.Cells.Item(13, 15).Value = 0.001555413647742459
$ws.Cells.Item(13, 16).Value = 0.001555413647742458
$ws.Cells.Item(13, 17).Value = 0.1057663592284444
$ws.Cells.Item(13, 18).Value = 0.951897233056
$ws.Cells.Item(13, 19).Value = 0.0003093848708750749
$ws.Cells.Item(13, 20).Value = 0.0003093848708750748
